$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clone the "J" attendance column into a new "K" column (new date) ---

# K1: new date header (same style/format as the other date headers)
$ws.Range("K1").Value = 45302
$ws.Range("K1").NumberFormat = $ws.Range("J1").NumberFormat

# K2:K6: attendance values for the new date column
$ws.Range("K2").Value = "Present"
$ws.Range("K3").Value = "Present"
$ws.Range("K4").Value = "Absent"
$ws.Range("K5").Value = "Absent"
$ws.Range("K6").Value = "Absent"

# --- Extend the existing "Present, Absent, Reason" list validation so it
#     also covers the new K column (H2:J6 -> H2:K6) ---
$ws.Range("C2:F6").Validation.Delete()
$ws.Range("H2:J6").Validation.Delete()
$ws.Range("C2:K6").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')
$ws.Range("G2:G6").Validation.Delete()

# --- Match the final selection left by the edit ---
$ws.Range("K6").Select()
